$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.714.45"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.07%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.732.14"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.51%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.15%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "242.43"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.78%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9984"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.4938"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.05%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2622"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06217"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.727.79"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.19%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.92"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.13%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07009"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6112"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.97%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.505"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.72%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "77.18"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9983"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.23%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.521.97"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.30%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.09%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.000007206"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.31%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.44"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.951.04"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.477"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.30%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.562"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.11%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.097"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.61%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "138.54"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.33%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.36"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.48%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.770"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.14%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.386"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "106.57"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.10%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.936"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.73%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07989"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.670"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.35%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.03%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 2)
$cell.NumberFormat = "@"
$cell.Value = "HuobiToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.609"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 2)
$cell.NumberFormat = "@"
$cell.Value = "ARBITRUM"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.58%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 2)
$cell.NumberFormat = "@"
$cell.Value = "ImmutableX"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6236"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 2)
$cell.NumberFormat = "@"
$cell.Value = "TrustWalletToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9422"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.71%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 2)
$cell.NumberFormat = "@"
$cell.Value = "RenderToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.038"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.59%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 2)
$cell.NumberFormat = "@"
$cell.Value = "MXToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.422"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.93%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 2)
$cell.NumberFormat = "@"
$cell.Value = "PaxDollar"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9984"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = "VeChain"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.01510"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.55%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = "FraxShare"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.578"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.10%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 2)
$cell.NumberFormat = "@"
$cell.Value = "Quant"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "99.49"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.74%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 2)
$cell.NumberFormat = "@"
$cell.Value = "TheSandbox"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.3860"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.21%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 2)
$cell.NumberFormat = "@"
$cell.Value = "Aptos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.939"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.77%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 2)
$cell.NumberFormat = "@"
$cell.Value = "Algorand"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1159"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.70%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = "@"
$cell.Value = "Cronos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05382"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.40%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = "@"
$cell.Value = "EnergySwap"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.891"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.10%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 2)
$cell.NumberFormat = "@"
$cell.Value = "Elrond"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "30.27"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.32%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 2)
$cell.NumberFormat = "@"
$cell.Value = "Aave"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "51.71"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.17%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 2)
$cell.NumberFormat = "@"
$cell.Value = "NEARProtocol"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.231"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.51%  "
$cell.Style = "Normal"
